$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")
$ws.Activate()

# --- Complete row 3 (finish the trip that was only half-filled) ---
$ws.Cells.Item(3,7).Value  = "12:48:57 PM"
$ws.Cells.Item(3,8).Value  = 1722440938028
$ws.Cells.Item(3,9).Value  = 1600
$ws.Cells.Item(3,10).Value = 1600
$ws.Cells.Item(3,11).Value = 0
$ws.Cells.Item(3,12).Value = "19 horas 51 minutos"

# --- Update row 4 (a subsequent edit to an already-completed trip) ---
$ws.Cells.Item(4,7).Value  = "1:05:38 PM"
$ws.Cells.Item(4,8).Value  = 1722441938606
$ws.Cells.Item(4,9).Value  = 1766
$ws.Cells.Item(4,10).Value = 1600
$ws.Cells.Item(4,11).Value = -166
$ws.Cells.Item(4,12).Value = "20 horas 7 minutos"

# Note: "2024-7-31" looks like a date to Excel's smart-entry, so plain
# assignment of that literal would be auto-converted into a date serial
# number. Copy it from an existing text cell that already holds it
# (e.g. F4) so the new cells stay plain text, matching the source data.
$dateCell = $ws.Cells.Item(4,6)

# --- New row 5: trip started but not finished ---
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "ABC1234"
$dateCell.Copy($ws.Cells.Item(5,3))
$ws.Cells.Item(5,4).Value = "1:04:49 PM"
$ws.Cells.Item(5,5).Value = 1722441889470

# --- New row 6: trip started but not finished ---
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "fsd"
$dateCell.Copy($ws.Cells.Item(6,3))
$ws.Cells.Item(6,4).Value = "1:22:10 PM"
$ws.Cells.Item(6,5).Value = 1722442931280

# --- New row 7: trip started and finished ---
$ws.Cells.Item(7,1).Value  = 6
$ws.Cells.Item(7,2).Value  = "fddf"
$dateCell.Copy($ws.Cells.Item(7,3))
$ws.Cells.Item(7,4).Value  = "1:25:12 PM"
$ws.Cells.Item(7,5).Value  = 1722443112893
$dateCell.Copy($ws.Cells.Item(7,6))
$ws.Cells.Item(7,7).Value  = "2:02:40 PM"
$ws.Cells.Item(7,8).Value  = 1722445360638
$ws.Cells.Item(7,9).Value  = 50
$ws.Cells.Item(7,10).Value = 40
$ws.Cells.Item(7,11).Value = -10
$ws.Cells.Item(7,12).Value = "0 horas 37 minutos"

# Leave the final active selection on F4, matching the saved view state
$ws.Range("F4").Select()
